$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "67.784.33", "0.140") that must
# stay as literal text rather than being auto-converted to a number by Excel,
# so those cells are pre-formatted as Text before the value is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.784.33'
$ws.Range("E2").Value = '  -3.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.824.57'
$ws.Range("E3").Value = '  -3.28%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.01'
$ws.Range("E5").Value = '  -2.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.56'
$ws.Range("E6").Value = '  -2.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.826.80'
$ws.Range("E7").Value = '  -3.08%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("E10").Value = '  -4.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.43'
$ws.Range("E11").Value = '  -0.75%  '

$ws.Range("E12").Value = '  -3.49%  '

$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.75'
$ws.Range("E14").Value = '  -4.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.465.48'
$ws.Range("E15").Value = '  -3.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.834.61'
$ws.Range("E16").Value = '  -2.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.975.03'
$ws.Range("E17").Value = '  -3.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.13'
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  -4.62%  '

$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("E21").Value = '  -2.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.68'
$ws.Range("E22").Value = '  -6.87%  '

$ws.Range("E23").Value = '  -2.72%  '

$ws.Range("E24").Value = '  -4.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.77'
$ws.Range("E25").Value = '  -4.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  -4.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.03'
$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.97'
$ws.Range("E29").Value = '  -3.46%  '

$ws.Range("E30").Value = '  -2.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.973.75'
$ws.Range("E31").Value = '  -3.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.56'
$ws.Range("E32").Value = '  -4.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.29'
$ws.Range("E33").Value = '  -6.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.94'
$ws.Range("E34").Value = '  -4.65%  '

$ws.Range("E35").Value = '  -1.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.785.30'
$ws.Range("E36").Value = '  -3.47%  '

$ws.Range("E37").Value = '  -4.99%  '

$ws.Range("E38").Value = '  +8.64%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.140'
$ws.Range("E39").Value = '  -1.05%  '

$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  -3.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.87'
$ws.Range("E41").Value = '  -5.20%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("E43").Value = '  -6.28%  '

$ws.Range("E44").Value = '  -7.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '417.84'
$ws.Range("E45").Value = '  -5.52%  '

$ws.Range("E46").Value = '  +5.70%  '

$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.60'
$ws.Range("E48").Value = '  -0.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '47.01'
$ws.Range("E49").Value = '  -2.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.14'
$ws.Range("E50").Value = '  -1.93%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0356'
$ws.Range("E51").Value = '  -3.65%  '
